$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = "AAAAA"
$ws.Range("B2").Value = "BBBB"
$ws.Range("C2").Value = "'12457845001"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "'12457855"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "Dr. Christian Santiago Izurieta Cruz"

# Update row 3
$ws.Range("A3").Value = "Msg. Jeremy Ismael Arias Benitez"
$ws.Range("B3").Value = "Alexander Francisco Tibanta Miranda"
$ws.Range("C3").Value = "'1728220441001"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = "'1728220441"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "Dr. Christian Santiago Izurieta Cruz"

# Delete rows 4-9 (remove extra data so dimension becomes A1:E3)
$ws.Rows("4:9").Delete()
